# Auto-generated Excel COM-interop script applying numeric updates
# to the Sargatanas_Profits market-price sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 11364138
$ws.Range("I41").Value = 20833670
$ws.Range("J41").Value = 700.8
$ws.Range("K41").Value = 20833670
$ws.Range("L41").Value = 700.8
$ws.Range("M41").Value = -20833230
$ws.Range("N41").Value = -1580.8
$ws.Range("H48").Value = 1222.5
$ws.Range("I48").Value = 1222.5
$ws.Range("K48").Value = 3667.5
$ws.Range("M48").Value = -3375.5
$ws.Range("H56").Value = 1222.5
$ws.Range("I56").Value = 1222.5
$ws.Range("K56").Value = 3667.5
$ws.Range("M56").Value = -3133.5
$ws.Range("H62").Value = 90927180
$ws.Range("I62").Value = 333334800
$ws.Range("J62").Value = 24333.625
$ws.Range("K62").Value = 333334800
$ws.Range("L62").Value = 24333.625
$ws.Range("M62").Value = -333334176
$ws.Range("N62").Value = -25581.625
$ws.Range("H65").Value = 90927180
$ws.Range("I65").Value = 333334800
$ws.Range("J65").Value = 24333.625
$ws.Range("K65").Value = 1666674000
$ws.Range("L65").Value = 121668.125
$ws.Range("M65").Value = -1666670880
$ws.Range("N65").Value = -127908.125
$ws.Range("H70").Value = 98487864
$ws.Range("I70").Value = 50002436
$ws.Range("K70").Value = 150007308
$ws.Range("M70").Value = -150007038
$ws.Range("H73").Value = 98487864
$ws.Range("I73").Value = 50002436
$ws.Range("K73").Value = 150007308
$ws.Range("M73").Value = -150006372
$ws.Range("H132").Value = 2001.6875
$ws.Range("I132").Value = 1927.963
$ws.Range("K132").Value = 5783.889
$ws.Range("M132").Value = -3253.889

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2206.07
$ws.Range("I32").Value = 2254.3125
$ws.Range("K32").Value = 2254.3125
$ws.Range("M32").Value = -1967.3125
$ws.Range("H102").Value = 804.53845
$ws.Range("I102").Value = 774.0909
$ws.Range("K102").Value = 774.0909
$ws.Range("M102").Value = 847.9091
$ws.Range("H110").Value = 55557308
$ws.Range("I110").Value = 1625
$ws.Range("K110").Value = 1625
$ws.Range("M110").Value = 420

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3951.2778
$ws.Range("I94").Value = 905.3333
$ws.Range("K94").Value = 905.3333
$ws.Range("M94").Value = -454.3333

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5138.5654
$ws.Range("I16").Value = 3582
$ws.Range("J16").Value = 6565.4165
$ws.Range("K16").Value = 3582
$ws.Range("L16").Value = 6565.4165
$ws.Range("M16").Value = -3295
$ws.Range("N16").Value = -7139.4165
$ws.Range("H58").Value = 9095774
$ws.Range("I58").Value = 14707470
$ws.Range("K58").Value = 14707470
$ws.Range("M58").Value = -14707267
$ws.Range("H99").Value = 9287.267
$ws.Range("I99").Value = 10549.667
$ws.Range("K99").Value = 10549.667
$ws.Range("M99").Value = -9051.666999999999
$ws.Range("H107").Value = 2843.4546
$ws.Range("I107").Value = 1439.6
$ws.Range("J107").Value = 4013.3333
$ws.Range("K107").Value = 1439.6
$ws.Range("L107").Value = 4013.3333
$ws.Range("N107").Value = -7853.3333
$ws.Range("M107").Value = 480.4000000000001
$ws.Range("H108").Value = 52020
$ws.Range("J108").Value = 52020
$ws.Range("L108").Value = 52020
$ws.Range("N108").Value = -59700
$ws.Range("H113").Value = 5138.5654
$ws.Range("I113").Value = 3582
$ws.Range("J113").Value = 6565.4165
$ws.Range("K113").Value = 3582
$ws.Range("L113").Value = 6565.4165
$ws.Range("M113").Value = -1412
$ws.Range("N113").Value = -10905.4165
$ws.Range("H126").Value = 9287.267
$ws.Range("I126").Value = 10549.667
$ws.Range("K126").Value = 31649.001
$ws.Range("M126").Value = -29179.001
$ws.Range("H132").Value = 5953.8975
$ws.Range("I132").Value = 2736.8
$ws.Range("J132").Value = 9340.315000000001
$ws.Range("K132").Value = 8210.400000000001
$ws.Range("L132").Value = 28020.945
$ws.Range("M132").Value = -5680.400000000001
$ws.Range("N132").Value = -33080.945
$ws.Range("H134").Value = 4633.25
$ws.Range("I134").Value = 2067.1282
$ws.Range("J134").Value = 8636.4
$ws.Range("K134").Value = 6201.3846
$ws.Range("L134").Value = 25909.2
$ws.Range("M134").Value = -3666.3846
$ws.Range("N134").Value = -30979.2
$ws.Range("H136").Value = 9095774
$ws.Range("I136").Value = 14707470
$ws.Range("K136").Value = 44122410
$ws.Range("M136").Value = -44119860

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 523.5
$ws.Range("I70").Value = 523.5
$ws.Range("K70").Value = 1570.5
$ws.Range("M70").Value = -1255.5
$ws.Range("H73").Value = 523.5
$ws.Range("I73").Value = 523.5
$ws.Range("K73").Value = 1570.5
$ws.Range("M73").Value = -478.5
$ws.Range("H87").Value = 1406.75
$ws.Range("I87").Value = 875.6667
$ws.Range("K87").Value = 2627.0001
$ws.Range("M87").Value = -1379.0001
$ws.Range("H90").Value = 1406.75
$ws.Range("I90").Value = 875.6667
$ws.Range("K90").Value = 7881.0003
$ws.Range("M90").Value = -1641.0003
$ws.Range("H107").Value = 50000396
$ws.Range("J107").Value = 100000350
$ws.Range("L107").Value = 300001050
$ws.Range("N107").Value = -300004890
$ws.Range("H113").Value = 1394.7222
$ws.Range("I113").Value = 1219.5
$ws.Range("K113").Value = 3658.5
$ws.Range("M113").Value = -1488.5
$ws.Range("H114").Value = 920.5833
$ws.Range("I114").Value = 852
$ws.Range("J114").Value = 943.44446
$ws.Range("K114").Value = 2556
$ws.Range("L114").Value = 2830.33338
$ws.Range("M114").Value = 698
$ws.Range("N114").Value = -9338.33338
$ws.Range("H132").Value = 8282.521000000001
$ws.Range("J132").Value = 13199.8
$ws.Range("L132").Value = 118798.2
$ws.Range("N132").Value = -123858.2
$ws.Range("H134").Value = 8038.727
$ws.Range("I134").Value = 7842.6
$ws.Range("K134").Value = 23527.8
$ws.Range("M134").Value = -18457.8

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2732.111
$ws.Range("J80").Value = 2977.7778
$ws.Range("L80").Value = 2977.7778
$ws.Range("N80").Value = -4973.7778
$ws.Range("H83").Value = 2732.111
$ws.Range("J83").Value = 2977.7778
$ws.Range("L83").Value = 14888.889
$ws.Range("N83").Value = -24872.889
$ws.Range("H132").Value = 4542.972
$ws.Range("I132").Value = 1400.1154
$ws.Range("K132").Value = 4200.3462
$ws.Range("M132").Value = -1670.3462
$ws.Range("H141").Value = 29452.545
$ws.Range("J141").Value = 34874.375
$ws.Range("L141").Value = 34874.375
$ws.Range("N141").Value = -45234.375

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1336.2
$ws.Range("I16").Value = 1336.2
$ws.Range("K16").Value = 1336.2
$ws.Range("M16").Value = -1166.2
$ws.Range("H61").Value = 4181.9614
$ws.Range("I61").Value = 1166.9
$ws.Range("J61").Value = 6066.375
$ws.Range("K61").Value = 1166.9
$ws.Range("L61").Value = 6066.375
$ws.Range("M61").Value = -964.9000000000001
$ws.Range("N61").Value = -6470.375
$ws.Range("H82").Value = 1085560.2
$ws.Range("I82").Value = 3521599.8
$ws.Range("J82").Value = 2876
$ws.Range("K82").Value = 3521599.8
$ws.Range("L82").Value = 2876
$ws.Range("M82").Value = -3521238.8
$ws.Range("N82").Value = -3598
$ws.Range("H85").Value = 1085560.2
$ws.Range("I85").Value = 3521599.8
$ws.Range("J85").Value = 2876
$ws.Range("K85").Value = 3521599.8
$ws.Range("L85").Value = 2876
$ws.Range("M85").Value = -3520351.8
$ws.Range("N85").Value = -5372
$ws.Range("H100").Value = 3588.5557
$ws.Range("I100").Value = 2185.4285
$ws.Range("K100").Value = 2185.4285
$ws.Range("M100").Value = -1644.4285
$ws.Range("H105").Value = 80000
$ws.Range("J105").Value = 80000
$ws.Range("L105").Value = 80000
$ws.Range("N105").Value = -86988
$ws.Range("H113").Value = 4181.9614
$ws.Range("I113").Value = 1166.9
$ws.Range("J113").Value = 6066.375
$ws.Range("K113").Value = 1166.9
$ws.Range("L113").Value = 6066.375
$ws.Range("M113").Value = 1003.1
$ws.Range("N113").Value = -10406.375
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents() | Out-Null
$ws.Range("H132").Value = 9621597
$ws.Range("I132").Value = 20002856
$ws.Range("J132").Value = 9319.777
$ws.Range("K132").Value = 60008568
$ws.Range("L132").Value = 27959.331
$ws.Range("M132").Value = -60006038
$ws.Range("N132").Value = -33019.331

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5726
$ws.Range("J74").Value = 5726
$ws.Range("L74").Value = 5726
$ws.Range("N74").Value = -7598
$ws.Range("H77").Value = 5726
$ws.Range("J77").Value = 5726
$ws.Range("L77").Value = 17178
$ws.Range("N77").Value = -26538
$ws.Range("H96").Value = 3000
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents() | Out-Null
$ws.Range("H122").Value = 3080.2163
$ws.Range("I122").Value = 1500.8334
$ws.Range("K122").Value = 4502.5002
$ws.Range("M122").Value = -2052.5002
$ws.Range("H123").Value = 45532.668
$ws.Range("J123").Value = 45532.668
$ws.Range("L123").Value = 45532.668
$ws.Range("N123").Value = -55332.668
$ws.Range("H132").Value = 10647609
$ws.Range("I132").Value = 12503730
$ws.Range("K132").Value = 37511190
$ws.Range("M132").Value = -37508660

